$d = $word.ActiveDocument

# --- Fill in the "Team Division of Labor" table rows 1-7 (row 0 is the header) ---
$t = $d.Tables.Item(1)

# Row 1: task / lead / collaborator
$t.Cell(2, 2).Range.Text = "Add escape option for each line when creating a new appliance "
$t.Cell(2, 3).Range.Text = "Destiny"
$t.Cell(2, 4).Range.Text = "Colin "

# Row 2: task / lead / collaborator
$t.Cell(3, 2).Range.Text = "Added print information of new appliance after it is created "
$t.Cell(3, 3).Range.Text = "Destiny"
$t.Cell(3, 4).Range.Text = "Colin "

# Row 3: task / lead / collaborator
$t.Cell(4, 2).Range.Text = "Added a method that deletes the newly created appliance if it was created incorrectly "
$t.Cell(4, 3).Range.Text = "Destiny"
$t.Cell(4, 4).Range.Text = "Colin "

# Row 4: task / lead only
$t.Cell(5, 2).Range.Text = "Added overloaded methods that print information to a csv file "
$t.Cell(5, 3).Range.Text = "Destiny"

# Row 5: task / lead only
$t.Cell(6, 2).Range.Text = "Implemented print to csv code that prints the user inputs, the program outputs, and other information to a csv file using overloaded methods. "
$t.Cell(6, 3).Range.Text = "Destiny"

# Row 6: task / lead only (lead has trailing space)
$t.Cell(7, 2).Range.Text = "Created program outline that details all inputs, outputs, and required methods and classes. "
$t.Cell(7, 3).Range.Text = "Destiny "

# Row 7: task / lead only (lead has trailing space)
$t.Cell(8, 2).Range.Text = "Created initial design slideshow and finished program outline and testing slides. "
$t.Cell(8, 3).Range.Text = "Destiny "

# --- Merge the "Tasks and responsibilities:" guidance sentence into a single run,
#     removing the gramStart/gramEnd proofErr markers around "E.g." ---
$d.Content.Find.Execute(
    "Be as specific as possible.  E.g. ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Be as specific as possible.  E.g. ",
    2) | Out-Null
